# Applies the cell-value updates described by the OOXML diff for cryptos.xlsx.
# Numeric-looking strings (e.g. "1.00", "7.48") are written via a temporary
# text NumberFormat so Excel keeps them as literal text instead of silently
# coercing them to numbers -- the original cells are all inline strings.
# The original Style is restored afterwards so no style/format cells change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.293.02"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "3.053.75"
$ws.Range("E3").Value = "  +2.90%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.05%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.64"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("E6").Value = "  +8.06%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.048.67"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("E9").Value = "  +0.86%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.39"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +10.68%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +7.48%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("E13").Value = "  +5.22%  "
$ws.Range("E14").Value = "  +4.18%  "
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").Value = "3.559.57"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "63.250.24"
$ws.Range("E17").Value = "  +3.29%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.07"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "3.046.69"
$ws.Range("E19").Value = "  +2.74%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.82"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("E21").Value = "  +3.34%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.697"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +2.24%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +2.79%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.89"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +2.06%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +6.06%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.91"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +10.03%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.32"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.48"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +9.51%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.73"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +3.68%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +5.83%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.79"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("E34").Value = "  +5.36%  "
$ws.Range("D35").Value = "0.0₃0865"
$ws.Range("E35").Value = "  +11.25%  "
$ws.Range("E36").Value = "  +3.00%  "
$ws.Range("E37").Value = "  +3.66%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +14.40%  "
$ws.Range("E39").Value = "  +3.42%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.65"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("E41").Value = "  +0.28%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.124"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +4.23%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.298"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +14.18%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.21"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +10.90%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "397.43"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("D47").Value = "2.754.91"
$ws.Range("E47").Value = "  +2.51%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.00"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("E50").Value = "  +3.73%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.24"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +4.30%  "
